# Cap nhat danh sach "San pham quan tam" khi import data mkt:
#  - Bo "Thao o to" / "Phu nano" / "Taplo o to" / "Thiet bi theo doi"
#  - Them "Tham 5D 6D" / "Tham TPE"

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("Sheet1"): cot G la "San pham quan tam" cho moi dong du lieu ---
$ws1.Range("G2").Value = "Thảm 5D 6D"
$ws1.Range("G3").Value = "Thảm 5D 6D"
$ws1.Range("G4").Value = "Thảm 5D 6D"
$ws1.Range("G5").Value = "Thảm 5D 6D"
$ws1.Range("G6").Value = "Thảm TPE"
$ws1.Range("G7").Value = "Thảm TPE"
$ws1.Range("G8").Value = "Thảm TPE"
$ws1.Range("G9").Value = "Thảm TPE"
$ws1.Range("G10").Value = "Thảm TPE"
# G10 trong ban goc lech style so voi cac o G khac; dong bo lai dinh dang
$ws1.Range("G9").Copy()
$ws1.Range("G10").PasteSpecial(-4122)  # xlPasteFormats

# --- Sheet2 ("Quy Tắc"): bang chu giai cac nhom san pham + xoa hyperlink cu ---
$ws2.Hyperlinks.Delete()

$ws2.Range("G3").Value = "Thảm 5D 6D"
$ws1.Range("G2").Copy()
$ws2.Range("G3").PasteSpecial(-4122)  # xlPasteFormats, dong bo voi style s=4

$ws2.Range("G4").Value = "Thảm TPE"

$ws2.Range("G5").Value = "(Có thể có thêm các Nhóm sản phẩm mới khi bổ sung trong cấu hình trên)"
$ws2.Range("G5").ClearFormats()

$ws2.Rows.Item(7).Delete()
$ws2.Rows.Item(6).Delete()

# --- Lua chon / tab dang active chuyen sang sheet "Quy Tắc" ---
$ws1.Range("G4").Select()
$ws2.Activate()
$ws2.Range("G16").Select()
